$d = $word.ActiveDocument
$xml = '<w:tbl><w:tblPr><w:tblW w:w="0" w:type="auto"/><w:tblBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tblBorders><w:tblCellMar><w:left w:w="0" w:type="dxa"/><w:right w:w="0" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/><w:tblCaption w:val=""/><w:tblDescription w:val=""/></w:tblPr><w:tblGrid><w:gridCol w:w="2742"/><w:gridCol w:w="1867"/><w:gridCol w:w="4397"/></w:tblGrid><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="1A027B75" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="76EF52D2" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Subject</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="564942A4" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Subtopic</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5472" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="538E6D43" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:bCs/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:u w:val="single"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Paper</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="5106A9D7" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="6C47025B" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="7B324881" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Mechanisms of Transmission</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5514" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="531468D6" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Wang et al., 2021, "Airborne Transmission of Respiratory Viruses"</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="5080805F" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="07AE7D9E" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="4C15C813" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Mechanisms of Transmission</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5658" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="7AA96B26" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Nogrady, 2024, "WHO redefines airborne transmission: what does that mean for future pandemics?"</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="28B40B9E" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="6CC9A82E" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="0F45ED97" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Mechanisms of Transmission</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5514" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="08B87D85" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Pohlker, 2021, "Respiratory aerosols and droplets in the transmission of infectious diseases"</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="5B720C90" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="1E2E7DFF" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="78A4720E" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Importance and Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5548" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="57B87BCB" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Poudel, 2021, "Impact of Covid-19 on health-related quality of life of patients: A structured review"</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="58928D3C" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="4158C6C8" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="05CEB1AA" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Importance and Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5514" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="6582736A" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Topcu, 2020, "The impact of COVID-19 on emerging stock markets"</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="009C0D86" w:rsidRPr="009C0D86" w14:paraId="3D596611" w14:textId="77777777" w:rsidTr="009C0D86"><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="09412DF0" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="248BFC53" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Importance and Impact</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5472" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p w14:paraId="295BAF1B" w14:textId="77777777" w:rsidR="009C0D86" w:rsidRPr="009C0D86" w:rsidRDefault="009C0D86" w:rsidP="009C0D86"><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r w:rsidRPr="009C0D86"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Dubey, 2020, "Psychosocial impact of COVID-19"</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Environmental Influences</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5662" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Pica, 2012, "Environmental factors affecting the transmission of respiratory viruses"</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3380" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Introduction to Airborne Disease Transmission Indoors</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2058" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Environmental Influences</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5514" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:left w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:bottom w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/><w:right w:val="single" w:sz="8" w:space="0" w:color="A3A3A3"/></w:tcBorders><w:tcMar><w:top w:w="40" w:type="dxa"/><w:left w:w="60" w:type="dxa"/><w:bottom w:w="40" w:type="dxa"/><w:right w:w="60" w:type="dxa"/></w:tcMar><w:hideMark/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/><w:kern w:val="0"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="en-GB"/><w14:ligatures w14:val="none"/></w:rPr><w:t>Ho, 2021, "Modeling airborne pathogen transport and transmission risks of SARS-CoV-2"</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="06618CA7" w14:textId="77777777" w:rsidR="006951C9" w:rsidRDefault="006951C9"/>'
$d.Content.InsertXML($xml)
